# Q3 Update - 2025
# UN-COI.xlsx: refresh the "fromCSV" sheet -
#  - new short-url token for the whole data pull (column B)
#  - row 16 now reports Dem. Rep. of the Congo (was Central African Rep.)
#  - row 17 now reports Sudan, a country new to this extract (was Dem. Rep. of the Congo)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- short-url refresh for every data row -----------------------------
$ws.Range("B2:B18").Value = "Mn0RyZ"

# --- row 16: Central African Rep. -> Dem. Rep. of the Congo -----------
# coo_id / coo_name / coo / coo_iso plus the refugees/returned_refugees split.
# These "numeric" id cells are stored as text in the workbook, so force a
# text number format before writing them (keeps them from being coerced to
# real numbers, matching the source data's column typing).
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "41"
$ws.Range("G16").Value = "Dem. Rep. of the Congo"
$ws.Range("H16").Value = "COD"
$ws.Range("I16").Value = "COD"
$ws.Range("O16").NumberFormat = "@"
$ws.Range("O16").Value = "5"
$ws.Range("P16").NumberFormat = "@"
$ws.Range("P16").Value = "0"

# --- row 17: Dem. Rep. of the Congo -> Sudan ---------------------------
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "177"
$ws.Range("G17").Value = "Sudan"
$ws.Range("H17").Value = "SUD"
$ws.Range("I17").Value = "SDN"
$ws.Range("O17").NumberFormat = "@"
$ws.Range("O17").Value = "8"
